$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Blad1")

# Row 41 ("Collapsing data") gains Swedish/English/Sitemap coverage marks
$ws.Range("C41").Value = 1
$ws.Range("D41").Value = 1
$ws.Range("E41").Value = 1

# Move the active selection to E41 to match the saved view state
$ws.Range("E41").Select()
